# This workbook holds a weekly price report. A new week of data (date
# 44607) was added to the top of the data block (row 44), pushing the
# existing rows down by two (old row 44 -> new row 46, etc.), and the
# oldest week that fell off the bottom (date 44595, old rows 162/163) was
# re-added at the very end of the table (new rows 164/165).
#
# Net effect: insert two new rows at row 44, shifting everything below
# down by two rows, then populate the two new rows with the new week's
# data (same market/category/etc. values as the row that is now directly
# below them, just with the new date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 44; everything currently at row 44 and
# below shifts down to row 46 and below (formatting/styles shift with
# the cells, same as using the Excel UI "Insert" command).
$ws.Rows("44:45").Insert()

# Fill in the new week's two rows (Primera / Segunda quality grades),
# matching the data pattern used throughout this sheet.
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 44607
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112040
$ws.Range("G44").Value = "Cilantro"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 200
$ws.Range("K44").Value = 600
$ws.Range("L44").Value = 700
$ws.Range("M44").Value = 650
$ws.Range("N44").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O44").Value = "Región de Ñuble"
$ws.Range("P44").Value = 650
$ws.Range("Q44").Value = 1
$ws.Range("R44").Value = "Hortaliza"

$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44607
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112040
$ws.Range("G45").Value = "Cilantro"
$ws.Range("H45").Value = "Sin especificar"
$ws.Range("I45").Value = "Segunda"
$ws.Range("J45").Value = 100
$ws.Range("K45").Value = 500
$ws.Range("L45").Value = 500
$ws.Range("M45").Value = 500
$ws.Range("N45").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O45").Value = "Región de Ñuble"
$ws.Range("P45").Value = 500
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = "Hortaliza"
